$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so that
# values such as "585.93" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.106.50"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "3.318.71"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "585.93"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("D6").Value = "181.97"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("D7").Value = "0.648"
$ws.Range("E7").Value = "  +2.52%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "3.319.19"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").Value = "0.126"
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("D11").Value = "6.79"
$ws.Range("E11").Value = "  +2.57%  "
$ws.Range("D12").Value = "0.404"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "3.895.18"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "0.131"
$ws.Range("E14").Value = "  -2.83%  "
$ws.Range("D15").Value = "66.148.75"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "26.27"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").Value = "3.336.84"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "423.96"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("D21").Value = "13.13"
$ws.Range("E21").Value = "  -2.78%  "
$ws.Range("D22").Value = "7.37"
$ws.Range("E22").Value = "  -2.49%  "
$ws.Range("D23").Value = "71.71"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").Value = "3.463.45"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").Value = "0.514"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").Value = "0.201"
$ws.Range("E28").Value = "  +4.73%  "
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("D30").Value = "8.93"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "1.91"
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("D33").Value = "22.36"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("E37").Value = "  -3.46%  "
$ws.Range("D38").Value = "160.81"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("D40").Value = "2.873.28"
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("D42").Value = "26.24"
$ws.Range("E42").Value = "  -4.25%  "
$ws.Range("D43").Value = "0.758"
$ws.Range("E43").Value = "  -4.50%  "
$ws.Range("D44").Value = "4.31"
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("D45").Value = "39.77"
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("E47").Value = "  -4.18%  "
$ws.Range("D48").Value = "2.30"
$ws.Range("E48").Value = "  -0.68%  "
$ws.Range("D49").Value = "23.19"
$ws.Range("E49").Value = "  -4.30%  "
$ws.Range("D50").Value = "313.82"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").Value = "0.0272"
$ws.Range("E51").Value = "  +0.03%  "
